# Auto-generated edit script: update crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the value to be stored as text even when it looks like a
    # plain number (e.g. "1.00", "12.90"), then strip the temporary
    # text NumberFormat so the cell style matches the untouched original
    # (no explicit numFmt / style index left behind).
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "53.272.04"
$ws.Range("E2").Value = "  +3.75%  "
$ws.Range("D3").Value = "3.153.06"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue "D5" "397.44"
$ws.Range("E5").Value = "  +3.41%  "
Set-TextValue "D6" "106.45"
$ws.Range("E6").Value = "  +3.99%  "
Set-TextValue "D7" "0.545"
$ws.Range("E7").Value = "  +0.36%  "
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +4.63%  "
Set-TextValue "D10" "38.78"
$ws.Range("E10").Value = "  +5.50%  "
$ws.Range("E11").Value = "  +1.22%  "
Set-TextValue "D12" "0.0871"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "3.649.36"
$ws.Range("E13").Value = "  +3.59%  "
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("E16").Value = "  +8.89%  "
$ws.Range("D17").Value = "3.156.12"
$ws.Range("E17").Value = "  +3.91%  "
Set-TextValue "D18" "10.79"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("D19").Value = "53.184.30"
$ws.Range("E19").Value = "  +3.41%  "
Set-TextValue "D20" "3.28"
$ws.Range("E20").Value = "  +4.51%  "
Set-TextValue "D21" "12.90"
$ws.Range("E21").Value = "  +4.03%  "
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  +1.27%  "
Set-TextValue "D23" "71.21"
$ws.Range("E23").Value = "  +1.47%  "
Set-TextValue "D24" "271.06"
$ws.Range("E24").Value = "  +1.23%  "
Set-TextValue "D25" "3.23"
$ws.Range("E25").Value = "  +2.46%  "
Set-TextValue "D26" "8.08"
$ws.Range("E26").Value = "  -1.31%  "
Set-TextValue "D27" "27.69"
$ws.Range("E27").Value = "  +3.27%  "
Set-TextValue "D28" "7.52"
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("E32").Value = "  +7.69%  "
Set-TextValue "D33" "37.41"
$ws.Range("E33").Value = "  +7.73%  "
Set-TextValue "D34" "0.0493"
$ws.Range("E34").Value = "  +10.78%  "
$ws.Range("E35").Value = "  +0.59%  "
Set-TextValue "D36" "50.50"
$ws.Range("E36").Value = "  +0.53%  "
Set-TextValue "D37" "3.53"
$ws.Range("E37").Value = "  +6.01%  "
Set-TextValue "D38" "0.999"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "2.76"
$ws.Range("E39").Value = "  +9.03%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D40" "4.15"
$ws.Range("E40").Value = "  +11.49%  "
$ws.Range("E41").Value = "  +2.12%  "
Set-TextValue "D42" "17.34"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("E43").Value = "  +1.64%  "
Set-TextValue "D44" "129.80"
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("E45").Value = "  +0.94%  "
Set-TextValue "D46" "22.37"
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("E47").Value = "  +0.09%  "
Set-TextValue "D48" "2.43"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").Value = "2.090.27"
$ws.Range("E49").Value = "  +3.28%  "
Set-TextValue "D50" "0.0506"
$ws.Range("E50").Value = "  +26.48%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-TextValue "D51" "0.0331"
$ws.Range("E51").Value = "  +4.68%  "
